$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values on specific rows
$ws.Range("F2").Value = -1
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -4
$ws.Range("F7").Value = -7
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -2
